# deep sea double count fix
# Updates computed landings/percentage values on rows 4, 5, and 7
# (Tier 1, Tier 2, and Global) to correct a double-counting bug.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Tier 1)
$ws.Range("B4").Value = 13.41533138794797
$ws.Range("C4").Value = 23.48790230289491
$ws.Range("D4").Value = 9.268718574650407
$ws.Range("E4").Value = 36.90323369084288
$ws.Range("F4").Value = 9.268718574650407
$ws.Range("G4").Value = 29.05515303058552
$ws.Range("H4").Value = 50.87049853954008
$ws.Range("I4").Value = 20.07434842987439
$ws.Range("J4").Value = 79.92565157012559
$ws.Range("K4").Value = 20.07434842987439

# Row 5 (Tier 2)
$ws.Range("C5").Value = 7.571422503565985
$ws.Range("D5").Value = 4.766509757150653
$ws.Range("E5").Value = 10.69128708186119
$ws.Range("F5").Value = 4.766509757150653
$ws.Range("G5").Value = 20.18311283805593
$ws.Range("H5").Value = 48.98125251884212
$ws.Range("I5").Value = 30.83563464310195
$ws.Range("J5").Value = 69.16436535689806
$ws.Range("K5").Value = 30.83563464310195

# Row 7 (Global)
$ws.Range("B7").Value = 18.50865755204852
$ws.Range("C7").Value = 33.60113412567404
$ws.Range("D7").Value = 17.34828565888902
$ws.Range("E7").Value = 52.10979167772258
$ws.Range("F7").Value = 17.34828565888902
$ws.Range("G7").Value = 26.64723565892969
$ws.Range("H7").Value = 48.37613624522655
$ws.Range("I7").Value = 24.97662809584376
$ws.Range("J7").Value = 75.02337190415625
$ws.Range("K7").Value = 24.97662809584376
